$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocol")

$ws.Range("D1").Value = "ANAT_or_FUNC"
$ws.Range("D5").Value = "func"
$ws.Range("D6").Value = "func"
$ws.Range("D7").Value = "func"
$ws.Range("D8").Value = "func"
$ws.Range("H7").Value = "func"
$ws.Range("H7").ClearFormats()
$ws.Range("H2").Value = "use for ANAT_or_FUNC"

$ws.Range("D24").Select()
